# ============================================================
# Edit: Update Step1_Data raw signal distribution values (tire-type
# filtering reprocessing), then recompute the dependent sheets:
#   - Step2_Sj            : row-wise running cumulative sum of Step1_Data
#   - Step3_DataPts_0.5/.7/.8/.9 : first-exceeds-threshold stats derived
#                                   from Step2_Sj's cumulative values
# ============================================================

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Step1_Data")
$ws2 = $wb.Worksheets.Item("Step2_Sj")

# ------------------------------------------------------------------
# 1) Write the reprocessed Step1_Data values (columns D:AJ, rows 2:6)
# ------------------------------------------------------------------
$ws1.Cells.Item(2,4).Value2 = 0.3302499310136194
$ws1.Cells.Item(2,5).Value2 = 0.1669039253837691
$ws1.Cells.Item(2,6).Value2 = 0.1557810055898387
$ws1.Cells.Item(2,8).Value2 = 0.08442668820431816
$ws1.Cells.Item(2,10).Value2 = 0.1883348156117799
$ws1.Cells.Item(2,11).Value2 = 0.01045107765183455
$ws1.Cells.Item(2,12).Value2 = 0.007678390196927198
$ws1.Cells.Item(2,13).Value2 = 0.01622912513676921
$ws1.Cells.Item(2,14).Value2 = 0.01837533895006768
$ws1.Cells.Item(2,21).Value2 = 0.0194848632221797
$ws1.Cells.Item(2,23).Value2 = 0.002084839038896349
$ws1.Cells.Item(3,4).Value2 = 0.3716032150744039
$ws1.Cells.Item(3,5).Value2 = 0.1120828411820291
$ws1.Cells.Item(3,6).Value2 = 0.2284767577242373
$ws1.Cells.Item(3,8).Value2 = 0.09814763423703488
$ws1.Cells.Item(3,9).Value2 = 0.01186764823945862
$ws1.Cells.Item(3,10).Value2 = 0.04550257655858475
$ws1.Cells.Item(3,11).Value2 = 0.08018781914404788
$ws1.Cells.Item(3,12).Value2 = 0.0256757149465948
$ws1.Cells.Item(3,14).Value2 = 0.006978414549993373
$ws1.Cells.Item(3,15).Value2 = 0.004661910947744816
$ws1.Cells.Item(3,16).Value2 = 0.01017194478993194
$ws1.Cells.Item(3,21).Value2 = 0.002925365375275477
$ws1.Cells.Item(3,31).Value2 = 0.001718157230663112
$ws1.Cells.Item(4,4).Value2 = 0.4049335784125176
$ws1.Cells.Item(4,5).Value2 = 0.1079340005213218
$ws1.Cells.Item(4,6).Value2 = 0.2163027720570184
$ws1.Cells.Item(4,7).Value2 = 0.03118590938781405
$ws1.Cells.Item(4,8).Value2 = 0.06110479488036279
$ws1.Cells.Item(4,10).Value2 = 0.06543893651627987
$ws1.Cells.Item(4,11).Value2 = 0.05799595403787781
$ws1.Cells.Item(4,12).Value2 = 0.01550286680366057
$ws1.Cells.Item(4,14).Value2 = 0.002530597755671145
$ws1.Cells.Item(4,16).Value2 = 0.002476277645915934
$ws1.Cells.Item(4,21).Value2 = 0.03459431198155977
$ws1.Cells.Item(5,4).Value2 = 0.06014753881389298
$ws1.Cells.Item(5,5).Value2 = 0.1915309565432505
$ws1.Cells.Item(5,6).Value2 = 0.26625077613972
$ws1.Cells.Item(5,7).Value2 = 0.1053771039220696
$ws1.Cells.Item(5,8).Value2 = 0.03532971587881437
$ws1.Cells.Item(5,9).Value2 = 0.0224069648525802
$ws1.Cells.Item(5,10).Value2 = 0.005173921023381215
$ws1.Cells.Item(5,11).Value2 = 0.1127109794310434
$ws1.Cells.Item(5,12).Value2 = 0.05891543941350453
$ws1.Cells.Item(5,13).Value2 = 0.02193094678144049
$ws1.Cells.Item(5,15).Value2 = 0.04993603835221524
$ws1.Cells.Item(5,22).Value2 = 0.05272002345521698
$ws1.Cells.Item(5,24).Value2 = 0.01233469435441936
$ws1.Cells.Item(5,25).Value2 = 0.004434868664624168
$ws1.Cells.Item(5,35).Value2 = 0.0008000323738270586
$ws1.Cells.Item(6,5).Value2 = 0.2907730386157671
$ws1.Cells.Item(6,6).Value2 = 0.2039501137812078
$ws1.Cells.Item(6,7).Value2 = 0.1586686203233068
$ws1.Cells.Item(6,8).Value2 = 0.02531671387237915
$ws1.Cells.Item(6,9).Value2 = 0.04944221589138056
$ws1.Cells.Item(6,11).Value2 = 0.06940606279183588
$ws1.Cells.Item(6,12).Value2 = 0.09585265178217162
$ws1.Cells.Item(6,13).Value2 = 0.03780375597064969
$ws1.Cells.Item(6,15).Value2 = 0.03325913465644115
$ws1.Cells.Item(6,22).Value2 = 0.02849876634905324
$ws1.Cells.Item(6,32).Value2 = 0.007028925965806903

# ------------------------------------------------------------------
# 2) Recompute Step2_Sj as the running cumulative sum (row-wise) of
#    Step1_Data, across columns B:AJ (col 2..36), rows 2:6.
# ------------------------------------------------------------------
$firstCol = 2   # column B
$lastCol  = 36  # column AJ

for ($row = 2; $row -le 6; $row++) {
    $running = 0.0
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $src = $ws1.Cells.Item($row, $col).Value2
        if ($null -eq $src) { $src = 0 }
        $running = $running + $src
        $ws2.Cells.Item($row, $col).Value2 = $running
    }
}

# ------------------------------------------------------------------
# 3) Recompute each Step3_DataPts_<threshold> sheet from the refreshed
#    Step1_Data / Step2_Sj values:
#      C = First_Noticeable_Increase_Index  (1-based offset from col D,
#          i.e. D=1,E=2,... ; first column whose Step1_Data value >= 0.05)
#      D = Point_Exceeds_Index              (0-based column index from A;
#          first column whose Step2_Sj cumulative value >= threshold)
#      E = First_Noticeable_Increase_Cumulative_Value (unchanged, = 0)
#      F = Point_Exceeds_Cumulative_Value   (the cumulative value at D)
#      G = Pulse_Width                      (= D - C)
# ------------------------------------------------------------------
$thresholdSheets = @(
    @{ Name = "Step3_DataPts_0.5"; Threshold = 0.5 },
    @{ Name = "Step3_DataPts_0.7"; Threshold = 0.7 },
    @{ Name = "Step3_DataPts_0.8"; Threshold = 0.8 },
    @{ Name = "Step3_DataPts_0.9"; Threshold = 0.9 }
)

$noticeThreshold = 0.05

foreach ($entry in $thresholdSheets) {
    $ws3 = $wb.Worksheets.Item($entry.Name)
    $threshold = $entry.Threshold

    for ($row = 2; $row -le 6; $row++) {

        # --- First_Noticeable_Increase_Index (column C) ---
        $fniIndex = $null
        for ($col = $firstCol; $col -le $lastCol; $col++) {
            $v = $ws1.Cells.Item($row, $col).Value2
            if ($null -eq $v) { $v = 0 }
            if ($v -ge $noticeThreshold) {
                $fniIndex = $col - 3   # col D (=4) -> 1
                break
            }
        }

        # --- Point_Exceeds_Index / Cumulative_Value (columns D, F) ---
        $exceedsIndex = $null
        $exceedsValue = $null
        for ($col = $firstCol; $col -le $lastCol; $col++) {
            $v = $ws2.Cells.Item($row, $col).Value2
            if ($null -eq $v) { $v = 0 }
            if ($v -ge $threshold) {
                $exceedsIndex = $col - 1   # column A = 0
                $exceedsValue = $v
                break
            }
        }

        $ws3.Cells.Item($row, 3).Value2 = $fniIndex                      # C
        $ws3.Cells.Item($row, 4).Value2 = $exceedsIndex                  # D
        $ws3.Cells.Item($row, 6).Value2 = $exceedsValue                  # F
        $ws3.Cells.Item($row, 7).Value2 = ($exceedsIndex - $fniIndex)    # G
    }
}
